$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the DHIS2 demo instance rows: versions 2.37/2.36/2.35 -> 2.40/2.39/2.38
# (both the visible label in column A and the URL / hyperlink text in column B)
$ws.Range("A2").Value = "DHIS2 Demo 2.40"
$ws.Range("B2").Value = "https://play.dhis2.org/2.40/"

$ws.Range("A3").Value = "DHIS2 Demo 2.39"
$ws.Range("B3").Value = "https://play.dhis2.org/2.39/"

$ws.Range("A4").Value = "DHIS2 Demo 2.38"
$ws.Range("B4").Value = "https://play.dhis2.org/2.38/"

# Leave the original Hyperlinks collection (rId1/rId2/rId3) untouched so the
# workbook keeps its existing hyperlink/style structure intact.

# Match the saved selection/active cell state from the edit.
$ws.Range("B8").Select()

$wb.Save()
